$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("equat_2")
$ws.Columns("F").ColumnWidth = 15.9
